# Generate Report for handback
# Adds the "925ce3e3-1d33-4722-a26b-623a1299e28d" file's handback status as a
# new row (row 4) on all three worksheets: Overview, zh-cn, de-de.

$wb = $excel.ActiveWorkbook

$newId = "925ce3e3-1d33-4722-a26b-623a1299e28d"
$newMd = "$newId.md"
$zhXlf = "$newId.3f46556f883e8e5b9fbf7aac5ef9aebfe3036e25.zh-cn.xlf"
$deXlf = "$newId.3f46556f883e8e5b9fbf7aac5ef9aebfe3036e25.de-de.xlf"

$zhHandoffDateTime  = "2016-02-16 04:06:38"
$zhHandbackDateTime = "2016-02-16 04:07:39"
$deHandoffDateTime  = "2016-02-16 04:06:53"
$deHandbackDateTime = "2016-02-16 04:08:08"

$statusInSync = "Handed back: in sync with en-US"
$hoReason     = "Include"

# Same blue (FF6495ED) used for the existing hyperlink-styled cells, encoded
# as the BGR-ordered long that Font.Color expects.
$hyperlinkColor = 15570276

function Style-AsHyperlink($range) {
    $range.Font.Underline = $true
    $range.Font.Color = $hyperlinkColor
}

# ---------------------------------------------------------------------------
# Sheet "Overview": File Name | zh-cn | de-de
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A4").Value = $newMd
Style-AsHyperlink $wsOverview.Range("A4")
$wsOverview.Hyperlinks.Add($wsOverview.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/925ce3e31d334722a26b623a1299e28d000000000/e2e/$newMd", "", "", $newMd) | Out-Null

$wsOverview.Range("B4").Value = $statusInSync
$wsOverview.Range("C4").Value = $statusInSync

# ---------------------------------------------------------------------------
# Sheet "zh-cn": Source File Name | Status | Correspond Handoff File |
#                Correspond Handoff Datetime | Target File |
#                Correspond Handback File | Correspond Handback DateTime |
#                Handoff Reason | Dependency From
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A4").Value = $newMd
Style-AsHyperlink $wsZh.Range("A4")
$wsZh.Hyperlinks.Add($wsZh.Range("A4"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/925ce3e31d334722a26b623a1299e28d000000000/e2e/$newMd", "", "", $newMd) | Out-Null

$wsZh.Range("B4").Value = $statusInSync

$wsZh.Range("C4").Value = $zhXlf
Style-AsHyperlink $wsZh.Range("C4")
$wsZh.Hyperlinks.Add($wsZh.Range("C4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/925ce3e31d334722a26b623a1299e28d000000000/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/$zhXlf", "", "", $zhXlf) | Out-Null

$wsZh.Range("D4").Value = $zhHandoffDateTime
$wsZh.Range("D4").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsZh.Range("E4").Value = $newMd
Style-AsHyperlink $wsZh.Range("E4")
$wsZh.Hyperlinks.Add($wsZh.Range("E4"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/925ce3e31d334722a26b623a1299e28d000000000/e2e/$newMd", "", "", $newMd) | Out-Null

$wsZh.Range("F4").Value = $zhXlf
Style-AsHyperlink $wsZh.Range("F4")
$wsZh.Hyperlinks.Add($wsZh.Range("F4"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/925ce3e31d334722a26b623a1299e28d000000000/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/$zhXlf", "", "", $zhXlf) | Out-Null

$wsZh.Range("G4").Value = $zhHandbackDateTime

$wsZh.Range("H4").Value = $hoReason

# ---------------------------------------------------------------------------
# Sheet "de-de": same columns as zh-cn
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A4").Value = $newMd
Style-AsHyperlink $wsDe.Range("A4")
$wsDe.Hyperlinks.Add($wsDe.Range("A4"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/925ce3e31d334722a26b623a1299e28d000000000/e2e/$newMd", "", "", $newMd) | Out-Null

$wsDe.Range("B4").Value = $statusInSync

$wsDe.Range("C4").Value = $deXlf
Style-AsHyperlink $wsDe.Range("C4")
$wsDe.Hyperlinks.Add($wsDe.Range("C4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/925ce3e31d334722a26b623a1299e28d000000000/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/$deXlf", "", "", $deXlf) | Out-Null

$wsDe.Range("D4").Value = $deHandoffDateTime
$wsDe.Range("D4").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsDe.Range("E4").Value = $newMd
Style-AsHyperlink $wsDe.Range("E4")
$wsDe.Hyperlinks.Add($wsDe.Range("E4"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/925ce3e31d334722a26b623a1299e28d000000000/e2e/$newMd", "", "", $newMd) | Out-Null

$wsDe.Range("F4").Value = $deXlf
Style-AsHyperlink $wsDe.Range("F4")
$wsDe.Hyperlinks.Add($wsDe.Range("F4"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/925ce3e31d334722a26b623a1299e28d000000000/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/$deXlf", "", "", $deXlf) | Out-Null

$wsDe.Range("G4").Value = $deHandbackDateTime

$wsDe.Range("H4").Value = $hoReason

Write-Host "Applied handback row for $newId"
